# New crime data collected - update weekly CompStat figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: volume number and report week dates
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/16/2024  Through  12/22/2024"

# ---------------------------------------------------------------------------
# Helper cells used as format donors for the handful of cells whose type
# flips between "numeric" and "text placeholder" (the sheet uses shared text
# "0" / "***.*" in place of a number when the underlying value is blank/NA).
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

function Set-TextPlaceholder($addr, $donorAddr, $text) {
    $donor = $ws.Range($donorAddr)
    $donor.Copy()
    $dst = $ws.Range($addr)
    $dst.PasteSpecial($xlPasteFormats)
    $dst.NumberFormat = "@"
    $dst.Value = $text
    $donor.Copy()
    $dst.PasteSpecial($xlPasteFormats)
}

function Set-NumericFromText($addr, $donorAddr, $value) {
    $donor = $ws.Range($donorAddr)
    $donor.Copy()
    $dst = $ws.Range($addr)
    $dst.PasteSpecial($xlPasteFormats)
    $dst.Value = $value
}

# ---------------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------------
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("L15").Value = -27.272727272727

# ---------------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------------
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 55
$ws.Range("K16").Value = 37.5
$ws.Range("L16").Value = 10
$ws.Range("M16").Value = -14.0625
$ws.Range("N16").Value = -79.853479853479

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 33.333333333333
$ws.Range("I17").Value = 151
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 25.833333333333
$ws.Range("L17").Value = 43.809523809523
$ws.Range("M17").Value = 77.647058823529
$ws.Range("N17").Value = -28.436018957346

# ---------------------------------------------------------------------------
# Row 18 - C18 becomes the "0" placeholder, D18/E18 become real numbers,
# F18 keeps its old value of 1 (formerly held by C18).
# ---------------------------------------------------------------------------
Set-TextPlaceholder "C18" "F22" "0"
Set-NumericFromText "D18" "D17" 2
Set-NumericFromText "E18" "E17" -100
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = -66.666666666666
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = 23.809523809523
$ws.Range("M18").Value = -7.142857142857
$ws.Range("N18").Value = -87.5

# ---------------------------------------------------------------------------
# Row 19 - C19 becomes a real number (it used to be the "0" placeholder).
# ---------------------------------------------------------------------------
Set-NumericFromText "C19" "D19" 3
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 200
$ws.Range("G19").Value = 9
$ws.Range("H19").Value = -33.333333333333
$ws.Range("I19").Value = 152
$ws.Range("J19").Value = 185
$ws.Range("K19").Value = -17.837837837837
$ws.Range("L19").Value = -10.059171597633
$ws.Range("M19").Value = 34.513274336283
$ws.Range("N19").Value = -21.243523316062

# ---------------------------------------------------------------------------
# Row 20
# ---------------------------------------------------------------------------
$ws.Range("D20").Value = 4
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -83.333333333333
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = 20.454545454545
$ws.Range("L20").Value = -10.169491525423
$ws.Range("N20").Value = -82.622950819672

# ---------------------------------------------------------------------------
# Row 21 (bold "TOTAL" row - uses styles 17/18, unaffected by the helper fns)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 10
$ws.Range("F21").Value = 21
$ws.Range("G21").Value = 29
$ws.Range("H21").Value = -27.586206896551
$ws.Range("I21").Value = 471
$ws.Range("J21").Value = 436
$ws.Range("K21").Value = 8.027522935779
$ws.Range("L21").Value = 4.899777282850
$ws.Range("M21").Value = 34.571428571428
$ws.Range("N21").Value = -66.737288135593

# ---------------------------------------------------------------------------
# Row 22 - F22 becomes the "0" placeholder (copy its own old format first).
# ---------------------------------------------------------------------------
Set-TextPlaceholder "F22" "C22" "0"
$ws.Range("H22").Value = -100

# ---------------------------------------------------------------------------
# Row 23 - G23 becomes "0" placeholder, H23 becomes "***.*" placeholder.
# ---------------------------------------------------------------------------
Set-TextPlaceholder "G23" "F23" "0"
Set-TextPlaceholder "H23" "E22" "***.*"
$ws.Range("L23").Value = 0

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 13
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 36
$ws.Range("G24").Value = 42
$ws.Range("H24").Value = -14.285714285714
$ws.Range("I24").Value = 494
$ws.Range("J24").Value = 471
$ws.Range("K24").Value = 4.883227176220
$ws.Range("L24").Value = -5.544933078393
$ws.Range("M24").Value = 67.457627118644

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -57.142857142857
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 250
$ws.Range("J25").Value = 237
$ws.Range("K25").Value = 5.485232067510
$ws.Range("L25").Value = 16.822429906542

# ---------------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = 138.461538461538
$ws.Range("I26").Value = 270
$ws.Range("J26").Value = 215
$ws.Range("K26").Value = 25.581395348837
$ws.Range("L26").Value = 26.168224299065
$ws.Range("M26").Value = -18.918918918918

# ---------------------------------------------------------------------------
# Row 27
# ---------------------------------------------------------------------------
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("L27").Value = -21.428571428571

# ---------------------------------------------------------------------------
# Row 28
# ---------------------------------------------------------------------------
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 19
$ws.Range("K28").Value = -5
$ws.Range("L28").Value = -17.391304347826

# ---------------------------------------------------------------------------
# Row 33 - D33 becomes "0" placeholder, E33 becomes "***.*" placeholder.
# ---------------------------------------------------------------------------
Set-TextPlaceholder "D33" "C33" "0"
Set-TextPlaceholder "E33" "H23" "***.*"

# ---------------------------------------------------------------------------
# A new week's footnote block is inserted as a blank row above the old
# row 56, pushing the final two footnote rows down to 57 and 58 (and the
# sheet dimension/merged cells follow automatically).
# ---------------------------------------------------------------------------
$ws.Rows(56).Insert()
$ws.Range("A56").Clear()
